$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nl = [char]10

# ---- Row 2 / Row 3 content updates (order chosen to mirror the
# original authoring sequence so shared-string ids line up naturally) ----

# C2: diagnosis question moves from D2 to C2
$ws.Range("C2").Value = "Maybe it only uploads scripts called by the main run script?"

# A3 (new row): the new problem entry
$ws.Range("A3").Value = "When publishing, get error 00068 ""script XX contains broken project data source: <feature layer name>"""

# D2: becomes a hyperlink labeled "Geonet thread" (replaces old ESRI follow-up text)
$ws.Hyperlinks.Add($ws.Range("D2"), "https://community.esri.com/", "", "", "Geonet thread")

# B2: append the new 3/6/2022 follow-up note to the existing solution text
$ws.Range("B2").Value = "3/1/2022 - Manually copying over the scripts that didn't work onto the server computer seems to fix." + $nl + "3/6/2022 - Instead of running ""from module import function"", run ""import module.function"" or ""import module""--this will enable that script to automatically copy over, but the 'utils' folder still doesn't publish, even with an __init__.py in it"

# B3: new solution/workaround text
$ws.Range("B3").Value = "Generally, try to find and tweak the string that may be causing the issue--usually it's because GIS is trying to confirm if it's a data source, and if it finds it is invalue, it will through the error." + $nl + "Workaround: make one offending url string a 1-item list, then pluck it back out of the list once it's being used--that way, ArcGIS thinks it's a list an doesn't scrutinize it like a string."

# D3: another "Geonet thread" hyperlink
$ws.Hyperlinks.Add($ws.Range("D3"), "https://community.esri.com/", "", "", "Geonet thread")

# ---- Formatting ----
$ws.Range("A3:B3").WrapText = $true

# ---- Row heights ----
$ws.Rows.Item(2).RowHeight = 90
$ws.Rows.Item(3).RowHeight = 120

# ---- Selection ----
$ws.Range("C3").Select() | Out-Null
